# Update the "想去人数" (want-to-go count) figures in the "展览" and
# "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 328
        4  = 8467
        5  = 6164
        6  = 534
        10 = 323
        11 = 1127
        12 = 83
    }
    "全部类型" = @{
        2  = 328
        4  = 8467
        5  = 6164
        6  = 534
        10 = 323
        15 = 1127
        16 = 83
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
